$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.820.92"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.457.78"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'548.47"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").Value = "'147.16"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("D9").Value = "2.457.49"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'5.40"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("D14").Value = "'26.17"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "2.901.68"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").Value = "61.555.28"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "2.455.52"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("D19").Value = "'10.94"
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "'320.33"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'1.89"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("D25").Value = "'63.90"
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("D26").Value = "0.0₃0980"
$ws.Range("E26").Value = "  -8.57%  "
$ws.Range("D27").Value = "2.581.11"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  -6.74%  "
$ws.Range("D30").Value = "'532.35"
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("D31").Value = "'8.24"
$ws.Range("E31").Value = "  -5.60%  "
$ws.Range("D32").Value = "'7.74"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  -6.61%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'4.77"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "'18.23"
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("D41").Value = "'1.76"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'140.41"
$ws.Range("E42").Value = "  -7.49%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'40.42"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "'2.26"
$ws.Range("E45").Value = "  -5.88%  "
$ws.Range("D46").Value = "'143.76"
$ws.Range("E46").Value = "  -5.53%  "
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("D48").Value = "'21.69"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'0.0930"
$ws.Range("E51").Value = "  -3.50%  "
